# Updating and Refactor Code
# Applies the aging-rate-2022 update: adds "Sep 2022" column data (column V)
# to the quality-of-AR table, shifts the manual adjustment row (row 8) over
# by one column to accommodate the new month, recomputes the dependent
# percentage/summary cells, and updates the instructional header text.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update header instruction text (shared string) ---
$ws.Range("B2").Value() = "STEP #1 isi angka kualitas AR dari bulan Januari 2021  -  Sep 2022"

# --- Fill in Sep 2022 (column V) figures for the four quality rows ---
$ws.Range("V4").Value() = 1865514000.0
$ws.Range("V5").Value() = 526249500.0
$ws.Range("V6").Value() = 977393380.0
$ws.Range("V7").Value() = 1442604640.0

# --- Row 8 (manual adjustment row): shift existing values one column to
#     the right (D..V -> E..W) and populate the newly vacated A8/D8 cells,
#     then set the new trailing figures for W8/X8 ---
$ws.Range("A8").Value() = 0.0
$ws.Range("D8").Value() = 0.0
$ws.Range("E8").Value() = 21094500.0
$ws.Range("F8").Value() = 32067480.0
$ws.Range("G8").Value() = -12559000.0
$ws.Range("H8").Value() = 18040996.0
$ws.Range("I8").Value() = 45577214.0
$ws.Range("J8").Value() = 9105000.0
$ws.Range("K8").Value() = 1353000.0
$ws.Range("L8").Value() = -2411000.0
$ws.Range("M8").Value() = 635954988.0
$ws.Range("N8").Value() = 168231959.0
$ws.Range("O8").Value() = -917340375.0
$ws.Range("P8").Value() = 10844000.0
$ws.Range("Q8").Value() = -70878980.0
$ws.Range("R8").Value() = 44126000.0
$ws.Range("S8").Value() = 54577000.0
$ws.Range("T8").Value() = 40290000.0
$ws.Range("U8").Value() = 36209880.0
$ws.Range("V8").Value() = -43491500.0
$ws.Range("W8").Value() = 208225500.0
$ws.Range("X8").Value() = -1442604640.0

# --- Row 9 (column totals): Sep 2022 total ---
$ws.Range("V9").Value() = 4811761520.0

# --- Row 13 (Lancar %): recomputed after adding Sep 2022 ---
$ws.Range("U13").Value() = 28.7407305523
$ws.Range("V13").Value() = 0.0
$ws.Range("W13").Value() = 50.2508201614
$ws.Range("Z13").Value() = 50.1288459753

# --- Row 14 (Kurang Lancar %): recomputed after adding Sep 2022 ---
$ws.Range("Q14").Value() = 75.0769551684
$ws.Range("U14").Value() = 0.0
$ws.Range("V14").Value() = 0.0
$ws.Range("W14").Value() = 34.4919513894
$ws.Range("Z14").Value() = 34.4919513894

# --- Row 15 (Diragukan %): recomputed after adding Sep 2022 ---
$ws.Range("S15").Value() = 1.0
$ws.Range("U15").Value() = 0.0
$ws.Range("V15").Value() = 0.0
$ws.Range("W15").Value() = 11.0189692235
$ws.Range("Z15").Value() = 11.0189692235

# --- PD summary figures (D19/D20) update, D21 cleared ---
$ws.Range("D19").Value() = 1.90986157
$ws.Range("D20").Value() = 3.8006575082
$ws.Range("D21").Value() = ""
